# Rename wc_lang.core.CrossReference -> wc_lang.core.DatabaseReference:
# rename the "Cross references" worksheet tab to "Database references".
# (Renaming the sheet automatically updates the sheet-scoped defined name
# "_xlnm._FilterDatabase" that refers to it via 'Cross references'!...)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cross references")
$ws.Name = "Database references"

# Make the renamed sheet the active / selected tab (moves tabSelected +
# activeTab from the previously active "Compartments" sheet to this one).
$ws.Activate()
